$d = $word.ActiveDocument

$d.Content.Find.Execute("49+20=69", $true, $true, $false, $false, $false, $true, 1, $false, "77-8=69", 2) | Out-Null
$d.Content.Find.Execute("7+65=72", $true, $true, $false, $false, $false, $true, 1, $false, "3+38=41", 2) | Out-Null
$d.Content.Find.Execute("89-66=23", $true, $true, $false, $false, $false, $true, 1, $false, "50-32=18", 2) | Out-Null
$d.Content.Find.Execute("52+30=82", $true, $true, $false, $false, $false, $true, 1, $false, "16-11=5", 2) | Out-Null
$d.Content.Find.Execute("76-46=30", $true, $true, $false, $false, $false, $true, 1, $false, "59-49=10", 2) | Out-Null
$d.Content.Find.Execute("24+29=53", $true, $true, $false, $false, $false, $true, 1, $false, "24+28=52", 2) | Out-Null
$d.Content.Find.Execute("39+19=58", $true, $true, $false, $false, $false, $true, 1, $false, "25+5=30", 2) | Out-Null
$d.Content.Find.Execute("43-23=20", $true, $true, $false, $false, $false, $true, 1, $false, "69-28=41", 2) | Out-Null
$d.Content.Find.Execute("57+2=59", $true, $true, $false, $false, $false, $true, 1, $false, "63-17=46", 2) | Out-Null
$d.Content.Find.Execute("22+1=23", $true, $true, $false, $false, $false, $true, 1, $false, "66-13=53", 2) | Out-Null
$d.Content.Find.Execute("42+50=92", $true, $true, $false, $false, $false, $true, 1, $false, "12+36=48", 2) | Out-Null
$d.Content.Find.Execute("94-79=15", $true, $true, $false, $false, $false, $true, 1, $false, "11-3=8", 2) | Out-Null
$d.Content.Find.Execute("4+20=24", $true, $true, $false, $false, $false, $true, 1, $false, "31+32=63", 2) | Out-Null
$d.Content.Find.Execute("82+16=98", $true, $true, $false, $false, $false, $true, 1, $false, "61-40=21", 2) | Out-Null
$d.Content.Find.Execute("35+15=50", $true, $true, $false, $false, $false, $true, 1, $false, "89-48=41", 2) | Out-Null
$d.Content.Find.Execute("3+69=72", $true, $true, $false, $false, $false, $true, 1, $false, "50+5=55", 2) | Out-Null
$d.Content.Find.Execute("15+27=42", $true, $true, $false, $false, $false, $true, 1, $false, "35+53=88", 2) | Out-Null
$d.Content.Find.Execute("1+89=90", $true, $true, $false, $false, $false, $true, 1, $false, "77-22=55", 2) | Out-Null
$d.Content.Find.Execute("65-61=4", $true, $true, $false, $false, $false, $true, 1, $false, "29+55=84", 2) | Out-Null
$d.Content.Find.Execute("52+14=66", $true, $true, $false, $false, $false, $true, 1, $false, "66-11=55", 2) | Out-Null
$d.Content.Find.Execute("91-38=53", $true, $true, $false, $false, $false, $true, 1, $false, "79-62=17", 2) | Out-Null
$d.Content.Find.Execute("86+0=86", $true, $true, $false, $false, $false, $true, 1, $false, "70+11=81", 2) | Out-Null
$d.Content.Find.Execute("35-25=10", $true, $true, $false, $false, $false, $true, 1, $false, "96-30=66", 2) | Out-Null
$d.Content.Find.Execute("41+48=89", $true, $true, $false, $false, $false, $true, 1, $false, "70+2=72", 2) | Out-Null
$d.Content.Find.Execute("99-92=7", $true, $true, $false, $false, $false, $true, 1, $false, "1+10=11", 2) | Out-Null
$d.Content.Find.Execute("66+18=84", $true, $true, $false, $false, $false, $true, 1, $false, "98-11=87", 2) | Out-Null
$d.Content.Find.Execute("92-87=5", $true, $true, $false, $false, $false, $true, 1, $false, "88-60=28", 2) | Out-Null
$d.Content.Find.Execute("11-8=3", $true, $true, $false, $false, $false, $true, 1, $false, "51+8=59", 2) | Out-Null
$d.Content.Find.Execute("75-57=18", $true, $true, $false, $false, $false, $true, 1, $false, "75+14=89", 2) | Out-Null
$d.Content.Find.Execute("28-13=15", $true, $true, $false, $false, $false, $true, 1, $false, "71+6=77", 2) | Out-Null
$d.Content.Find.Execute("20+2=22", $true, $true, $false, $false, $false, $true, 1, $false, "7+10=17", 2) | Out-Null
$d.Content.Find.Execute("35+55=90", $true, $true, $false, $false, $false, $true, 1, $false, "54+31=85", 2) | Out-Null
$d.Content.Find.Execute("23+50=73", $true, $true, $false, $false, $false, $true, 1, $false, "14+68=82", 2) | Out-Null
$d.Content.Find.Execute("33+65=98", $true, $true, $false, $false, $false, $true, 1, $false, "7+66=73", 2) | Out-Null
$d.Content.Find.Execute("75-1=74", $true, $true, $false, $false, $false, $true, 1, $false, "32+62=94", 2) | Out-Null
$d.Content.Find.Execute("90-8=82", $true, $true, $false, $false, $false, $true, 1, $false, "99-73=26", 2) | Out-Null
$d.Content.Find.Execute("2+40=42", $true, $true, $false, $false, $false, $true, 1, $false, "45+34=79", 2) | Out-Null
$d.Content.Find.Execute("65+20=85", $true, $true, $false, $false, $false, $true, 1, $false, "64-48=16", 2) | Out-Null
$d.Content.Find.Execute("80-47=33", $true, $true, $false, $false, $false, $true, 1, $false, "52+47=99", 2) | Out-Null
$d.Content.Find.Execute("1+64=65", $true, $true, $false, $false, $false, $true, 1, $false, "14+1=15", 2) | Out-Null
$d.Content.Find.Execute("6+69=75", $true, $true, $false, $false, $false, $true, 1, $false, "50+32=82", 2) | Out-Null
$d.Content.Find.Execute("74+6=80", $true, $true, $false, $false, $false, $true, 1, $false, "18+45=63", 2) | Out-Null
$d.Content.Find.Execute("44-16=28", $true, $true, $false, $false, $false, $true, 1, $false, "66-63=3", 2) | Out-Null
$d.Content.Find.Execute("58-19=39", $true, $true, $false, $false, $false, $true, 1, $false, "66-39=27", 2) | Out-Null
$d.Content.Find.Execute("69+30=99", $true, $true, $false, $false, $false, $true, 1, $false, "71+25=96", 2) | Out-Null
$d.Content.Find.Execute("81+0=81", $true, $true, $false, $false, $false, $true, 1, $false, "80-36=44", 2) | Out-Null
$d.Content.Find.Execute("98-54=44", $true, $true, $false, $false, $false, $true, 1, $false, "36-25=11", 2) | Out-Null
$d.Content.Find.Execute("0+95=95", $true, $true, $false, $false, $false, $true, 1, $false, "96-13=83", 2) | Out-Null
$d.Content.Find.Execute("95-27=68", $true, $true, $false, $false, $false, $true, 1, $false, "40+51=91", 2) | Out-Null
$d.Content.Find.Execute("93-93=0", $true, $true, $false, $false, $false, $true, 1, $false, "65-10=55", 2) | Out-Null
$d.Content.Find.Execute("86-78=8", $true, $true, $false, $false, $false, $true, 1, $false, "79-38=41", 2) | Out-Null
$d.Content.Find.Execute("9+79=88", $true, $true, $false, $false, $false, $true, 1, $false, "12+28=40", 2) | Out-Null
$d.Content.Find.Execute("1+88=89", $true, $true, $false, $false, $false, $true, 1, $false, "13+71=84", 2) | Out-Null
$d.Content.Find.Execute("11+19=30", $true, $true, $false, $false, $false, $true, 1, $false, "17+37=54", 2) | Out-Null
$d.Content.Find.Execute("65+3=68", $true, $true, $false, $false, $false, $true, 1, $false, "46-13=33", 2) | Out-Null
$d.Content.Find.Execute("6+32=38", $true, $true, $false, $false, $false, $true, 1, $false, "44-11=33", 2) | Out-Null
$d.Content.Find.Execute("98-45=53", $true, $true, $false, $false, $false, $true, 1, $false, "74-28=46", 2) | Out-Null
$d.Content.Find.Execute("35+56=91", $true, $true, $false, $false, $false, $true, 1, $false, "13+12=25", 2) | Out-Null
$d.Content.Find.Execute("41+51=92", $true, $true, $false, $false, $false, $true, 1, $false, "3+95=98", 2) | Out-Null
$d.Content.Find.Execute("88-61=27", $true, $true, $false, $false, $false, $true, 1, $false, "22+27=49", 2) | Out-Null
$d.Content.Find.Execute("54-39=15", $true, $true, $false, $false, $false, $true, 1, $false, "39+54=93", 2) | Out-Null
$d.Content.Find.Execute("51-26=25", $true, $true, $false, $false, $false, $true, 1, $false, "61-59=2", 2) | Out-Null
$d.Content.Find.Execute("79-42=37", $true, $true, $false, $false, $false, $true, 1, $false, "33+6=39", 2) | Out-Null
$d.Content.Find.Execute("88-35=53", $true, $true, $false, $false, $false, $true, 1, $false, "57-27=30", 2) | Out-Null
$d.Content.Find.Execute("54+9=63", $true, $true, $false, $false, $false, $true, 1, $false, "20+46=66", 2) | Out-Null
$d.Content.Find.Execute("98-60=38", $true, $true, $false, $false, $false, $true, 1, $false, "86-10=76", 2) | Out-Null
$d.Content.Find.Execute("3-1=2", $true, $true, $false, $false, $false, $true, 1, $false, "57+20=77", 2) | Out-Null
$d.Content.Find.Execute("47+40=87", $true, $true, $false, $false, $false, $true, 1, $false, "33+54=87", 2) | Out-Null
$d.Content.Find.Execute("94-84=10", $true, $true, $false, $false, $false, $true, 1, $false, "83-59=24", 2) | Out-Null
$d.Content.Find.Execute("82-32=50", $true, $true, $false, $false, $false, $true, 1, $false, "77-71=6", 2) | Out-Null
$d.Content.Find.Execute("91-3=88", $true, $true, $false, $false, $false, $true, 1, $false, "67+18=85", 2) | Out-Null
$d.Content.Find.Execute("78-59=19", $true, $true, $false, $false, $false, $true, 1, $false, "91-88=3", 2) | Out-Null
$d.Content.Find.Execute("18-2=16", $true, $true, $false, $false, $false, $true, 1, $false, "8+42=50", 2) | Out-Null
$d.Content.Find.Execute("96-12=84", $true, $true, $false, $false, $false, $true, 1, $false, "17+62=79", 2) | Out-Null
$d.Content.Find.Execute("63-22=41", $true, $true, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$d.Content.Find.Execute("86-50=36", $true, $true, $false, $false, $false, $true, 1, $false, "4+59=63", 2) | Out-Null
$d.Content.Find.Execute("2+82=84", $true, $true, $false, $false, $false, $true, 1, $false, "64+16=80", 2) | Out-Null
$d.Content.Find.Execute("46-43=3", $true, $true, $false, $false, $false, $true, 1, $false, "71-15=56", 2) | Out-Null
$d.Content.Find.Execute("99-5=94", $true, $true, $false, $false, $false, $true, 1, $false, "2-0=2", 2) | Out-Null
$d.Content.Find.Execute("60-36=24", $true, $true, $false, $false, $false, $true, 1, $false, "19+69=88", 2) | Out-Null
$d.Content.Find.Execute("60-24=36", $true, $true, $false, $false, $false, $true, 1, $false, "41+22=63", 2) | Out-Null
$d.Content.Find.Execute("93-77=16", $true, $true, $false, $false, $false, $true, 1, $false, "91-40=51", 2) | Out-Null
$d.Content.Find.Execute("4+73=77", $true, $true, $false, $false, $false, $true, 1, $false, "67-5=62", 2) | Out-Null
$d.Content.Find.Execute("67-58=9", $true, $true, $false, $false, $false, $true, 1, $false, "9+89=98", 2) | Out-Null
$d.Content.Find.Execute("61-50=11", $true, $true, $false, $false, $false, $true, 1, $false, "71-56=15", 2) | Out-Null
$d.Content.Find.Execute("72-53=19", $true, $true, $false, $false, $false, $true, 1, $false, "17+81=98", 2) | Out-Null
$d.Content.Find.Execute("82-67=15", $true, $true, $false, $false, $false, $true, 1, $false, "64-4=60", 2) | Out-Null
$d.Content.Find.Execute("18-12=6", $true, $true, $false, $false, $false, $true, 1, $false, "8+23=31", 2) | Out-Null
$d.Content.Find.Execute("64-18=46", $true, $true, $false, $false, $false, $true, 1, $false, "62-24=38", 2) | Out-Null
$d.Content.Find.Execute("34-34=0", $true, $true, $false, $false, $false, $true, 1, $false, "23+67=90", 2) | Out-Null
$d.Content.Find.Execute("28+46=74", $true, $true, $false, $false, $false, $true, 1, $false, "44+48=92", 2) | Out-Null
$d.Content.Find.Execute("28+29=57", $true, $true, $false, $false, $false, $true, 1, $false, "92-73=19", 2) | Out-Null
$d.Content.Find.Execute("75+19=94", $true, $true, $false, $false, $false, $true, 1, $false, "19-5=14", 2) | Out-Null
$d.Content.Find.Execute("13-10=3", $true, $true, $false, $false, $false, $true, 1, $false, "4+6=10", 2) | Out-Null
$d.Content.Find.Execute("70-58=12", $true, $true, $false, $false, $false, $true, 1, $false, "34-19=15", 2) | Out-Null
$d.Content.Find.Execute("71-41=30", $true, $true, $false, $false, $false, $true, 1, $false, "84-44=40", 2) | Out-Null
$d.Content.Find.Execute("36+41=77", $true, $true, $false, $false, $false, $true, 1, $false, "45+0=45", 2) | Out-Null
$d.Content.Find.Execute("54-38=16", $true, $true, $false, $false, $false, $true, 1, $false, "30+54=84", 2) | Out-Null
$d.Content.Find.Execute("66-66=0", $true, $true, $false, $false, $false, $true, 1, $false, "17+14=31", 2) | Out-Null
$d.Content.Find.Execute("59-53=6", $true, $true, $false, $false, $false, $true, 1, $false, "57-52=5", 2) | Out-Null
